# Performance optimisation for missing forms in form matrix.
#
# The "sub-step" rows (3,4,5 and 7,8,9) had their Test Step / Expected
# Results text sitting in columns E/F instead of A/B (with B..G otherwise
# blank). Move that text into A/B and drop the now-unused C:G cells on
# those rows. Then append a third scenario (rows 10-13) using the same
# template, and resize a couple of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-SubRow($row, $step, $expected) {
    $ws.Range("A$row").Value = $step
    $ws.Range("B$row").Value = $expected
    $ws.Range("C$row").Value = $null
    $ws.Range("D$row").Value = $null
    $ws.Range("E$row").Value = $null
    $ws.Range("F$row").Value = $null
    $ws.Range("G$row").Value = $null
}

Fix-SubRow 3 "initiate a New Submission transaction for a CA policy" "User should be able to start a New submission transaction"
Fix-SubRow 4 "Add  will triggered when Manhole Liability coverage is selected" "User should be able to add all the terms"
Fix-SubRow 5 "Verify that the below Coverage Term(s) is(are) displayed FIELDS" "The Coverage Term(s) should be displayed along with the default value (if any) and options available for selection"

Fix-SubRow 7 "initiate a New Submission transaction for a CA policy" "User should be able to start a New submission transaction"
Fix-SubRow 8 "Add  will triggered when Manhole Liability coverage is selected" "User should be able to add all the terms"
Fix-SubRow 9 "Verify that the below Coverage Term(s) is(are) displayed FIELDS" "The Coverage Term(s) should be displayed along with the default value (if any) and options available for selection"

# New scenario block (rows 10-13), mirroring the GA/AZ blocks above.
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "TC_1_Verify the dynamic fields for BAS BAS EDITION - FORMNAME"
$ws.Range("C10").Value = "Verify the dynamic fields for BAS BAS EDITION - FORMNAME"
$ws.Range("D10").Value = "Manual"
$ws.Range("E10").Value = "Login to PC and initiate a submission for STATE"
$ws.Range("F10").Value = "User should be able to log in successfully and should be navigated to the Home Screen"
$ws.Range("G10").Value = "SmartComm/Dyanamic Forms/SBGWI-256 - Something"

Fix-SubRow 11 "initiate a New Submission transaction for a CA policy" "User should be able to start a New submission transaction"
Fix-SubRow 12 "Add TRIGGERING_CONDITION" "User should be able to add all the terms"
Fix-SubRow 13 "Verify that the below Coverage Term(s) is(are) displayed FIELDS" "The Coverage Term(s) should be displayed along with the default value (if any) and options available for selection"

# Column width adjustments: A widens to match B/C (50), E narrows to 47.
$ws.Columns.Item(1).ColumnWidth = 49.2
$ws.Columns.Item(5).ColumnWidth = 46.2
